$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Text rename sweep: "Web Data Sanity" -> "Sanity Suite Test"
#    (every cell sharing these strings gets rewritten so the shared
#    string table collapses back down to the same text everywhere,
#    exactly like the source diff shows for sharedStrings.xml)
# ------------------------------------------------------------------
$wsCurrency  = $wb.Worksheets.Item("AddCurrency")
$wsCategory  = $wb.Worksheets.Item("AddProductCategory1")
$wsCustomer  = $wb.Worksheets.Item("AddCustomer")

$wsCurrency.Range("C1").Value = "Sanity Suite Test"
$wsCurrency.Range("D1").Value = "Sanity Suite Test Child"
$wsCurrency.Range("F1").Value = "Working as admin Sanity Suite Test Child X"

$wsCategory.Range("C1").Value = "Sanity Suite Test"
$wsCategory.Range("E1").Value = "Sanity Suite Test Child"

$wsCustomer.Range("C1").Value = "Sanity Suite Test"
$wsCustomer.Range("C2").Value = "Sanity Suite Test Child"

# "Canadian Dollar " (trailing space) -> "Canadian Dollar"
$wsCustomer.Range("F2").Value = "Canadian Dollar"

# ------------------------------------------------------------------
# 2) New sheet: CreateOrder (4th tab, becomes the active tab)
# ------------------------------------------------------------------
$wsOrder = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsOrder.Name = "CreateOrder"

# Row 1
$wsOrder.Range("A1").Value = "admin"
$wsOrder.Range("B1").Value = "WebData@123"
$wsOrder.Range("C1").Value = "Sanity Suite Test"
$wsOrder.Range("D1").Value = "USD Customer"
$wsOrder.Range("E1").Value = "''01/01/1970'"
$wsOrder.Range("F1").Value = "Test Product1"
$wsOrder.Range("G1").Value = "Test Product1"
$wsOrder.Range("H1").Value = 90
$wsOrder.Range("I1").Value = 120
$wsOrder.Range("J1").Value = 93
$wsOrder.Range("K1").Value = 123
$wsOrder.Range("L1").Value = 96
$wsOrder.Range("M1").Value = 126

# Row 2
$wsOrder.Range("C2").Value = "Sanity Suite Test Child"
$wsOrder.Range("D2").Value = "CAD Customer"
$wsOrder.Range("E2").Value = "''01/26/2017'"
$wsOrder.Range("F2").Value = "Test Product2"
$wsOrder.Range("G2").Value = "Test Product3"
$wsOrder.Range("H2").Value = 90
$wsOrder.Range("I2").Value = 120
$wsOrder.Range("J2").Value = 93
$wsOrder.Range("K2").Value = 123
$wsOrder.Range("L2").Value = 96
$wsOrder.Range("M2").Value = 126

# Row 3
$wsOrder.Range("D3").Value = "Euro Customer"
$wsOrder.Range("E3").Value = "''06/26/2018'"
$wsOrder.Range("F3").Value = "Test Product4"
$wsOrder.Range("G3").Value = "Test Product4"
$wsOrder.Range("H3").Value = 99
$wsOrder.Range("I3").Value = 129
$wsOrder.Range("J3").Value = 102
$wsOrder.Range("K3").Value = 132
$wsOrder.Range("L3").Value = 105
$wsOrder.Range("M3").Value = 135

# Row 4
$wsOrder.Range("E4").Value = "''01/01/2020'"
$wsOrder.Range("H4").Value = 108
$wsOrder.Range("I4").Value = 138
$wsOrder.Range("J4").Value = 111
$wsOrder.Range("K4").Value = 141
$wsOrder.Range("L4").Value = 114
$wsOrder.Range("M4").Value = 144

# Hyperlink on B1, matching the other sheets' mailto link pattern
$wsOrder.Hyperlinks.Add($wsOrder.Range("B1"), "mailto:WebData@123", $null, $null, "WebData@123")

# ------------------------------------------------------------------
# 3) Cursor / selection bookkeeping (best-effort; matches final
#    active cell reported in the diff for each sheet)
# ------------------------------------------------------------------
$wsCurrency.Range("F1").Select()
$wsCategory.Range("G9").Select()
$wsCustomer.Range("E1").Select()
$wsOrder.Range("E10").Select()
$wsOrder.Activate()
